$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range('C2').Value = '[-, -, ''MEC-3B-Mec. Manut.Equip. ind.'', -]'
$ws.Range('D2').Value = '-'
$ws.Range('E2').Value = '[-, -, ''MEC-3B-M.S.R. ar Cond.'', -]'
$ws.Range('E3').Value = '[-, -, ''MEC-3B-M.S.R. ar Cond.'', -]'
$ws.Range('F3').Value = '-'
$ws.Range('E4').Value = '[-, -, ''MEC-3B-M.S.R. ar Cond.'', -]'
$ws.Range('F4').Value = '-'
$ws.Range('E6').Value = '[-, -, -, ''MEC-3B-Mec. Manut.Equip. ind.'']'
$ws.Range('F6').Value = '-'
$ws.Range('E7').Value = '[-, -, -, ''MEC-3B-Mec. Manut.Equip. ind.'']'
$ws.Range('C8').Value = '[''MEC-3B-Mec. Manut.Equip. ind.'', -, -, -]'
$ws.Range('E8').Value = '[-, ''MEC-3B-M.S.R. ar Cond.'', -, -]'
$ws.Range('B18').Value = '-'
$ws.Range('C18').Value = '[-, ''MEC-2NA-M.S.R.A.C.'', -, -]'
$ws.Range('D18').Value = '[Ismail-Metrologia 2-2NB, ''MEC-2NA-M.S.R.A.C.'', ''MEC-2NB-M.S.R.A.C.'', -]'
$ws.Range('E18').Value = '[''MEC-2NB-M.S.R.A.C.'', ''ELM-2NA-Sistemas de Refrigeração'', ''MEC-2NB-M. Maq. E. I.'', -]'
$ws.Range('F18').Value = '-'
$ws.Range('B19').Value = '-'
$ws.Range('C19').Value = '[-, ''MEC-2NA-M.S.R.A.C.'', -, -]'
$ws.Range('D19').Value = '[''MEC-2NB-M. Maq. E. I.'', ''MEC-2NA-M.S.R.A.C.'', Leandro-M.S.R.A.C.-2NB, -]'
$ws.Range('E19').Value = '[-, ''ELM-2NA-Sistemas de Refrigeração'', -, -]'
$ws.Range('C20').Value = '[-, ''ELM-2NA-Sistemas de Refrigeração'', -, -]'
$ws.Range('D20').Value = '[''MEC-2NB-M. Maq. E. I.'', -, ''MEC-2NB-M.S.R.A.C.'', -]'
$ws.Range('F20').Value = '-'
$ws.Range('C21').Value = '[-, ''ELM-2NA-Sistemas de Refrigeração'', -, -]'
$ws.Range('D21').Value = '-'
$ws.Range('E21').Value = '-'
$ws.Range('F21').Value = '-'
